$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the old row 228, shifting the existing
# rows 228-233 down to 230-235. Excel's Insert() copies formatting from the
# row above, which matches what the target workbook expects (date style
# preserved in column D).
$ws.Rows.Item(228).Insert()
$ws.Rows.Item(228).Insert()

# --- New row 228: Ciboulette, Primera, week of 2021-09-09 ---
$ws.Range("A228").Value = 9
$ws.Range("B228").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C228").Value = "Metropolitana"
$ws.Range("D228").Value = 44448
$ws.Range("E228").Value = 13
$ws.Range("F228").Value = 100112039
$ws.Range("G228").Value = "Ciboulette"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 250
$ws.Range("K228").Value = 1500
$ws.Range("L228").Value = 2000
$ws.Range("M228").Value = 1750
$ws.Range("N228").Value = "`$/docena de atados"
$ws.Range("O228").Value = "Región Metropolitana"
$ws.Range("P228").Value = 583
$ws.Range("Q228").Value = 3
$ws.Range("R228").Value = "Hortaliza"

# --- New row 229: Ciboulette, Segunda, week of 2021-09-09 ---
$ws.Range("A229").Value = 9
$ws.Range("B229").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C229").Value = "Metropolitana"
$ws.Range("D229").Value = 44448
$ws.Range("E229").Value = 13
$ws.Range("F229").Value = 100112039
$ws.Range("G229").Value = "Ciboulette"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Segunda"
$ws.Range("J229").Value = 133
$ws.Range("K229").Value = 1100
$ws.Range("L229").Value = 1400
$ws.Range("M229").Value = 1251
$ws.Range("N229").Value = "`$/docena de atados"
$ws.Range("O229").Value = "Región Metropolitana"
$ws.Range("P229").Value = 417
$ws.Range("Q229").Value = 3
$ws.Range("R229").Value = "Hortaliza"
